$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.237.27'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").Value = '3.148.73'
$ws.Range("E3").Value = '  +2.16%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.82'
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '618.69'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("E7").Value = '  +4.98%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.368'
$ws.Range("E8").Value = '  +2.80%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = '3.144.73'
$ws.Range("E10").Value = '  +1.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.734'
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.203'
$ws.Range("E12").Value = '  +3.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000245'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.09'
$ws.Range("E14").Value = '  +0.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.52'
$ws.Range("E15").Value = '  +2.57%  '
$ws.Range("D16").Value = '90.442.33'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").Value = '3.752.06'
$ws.Range("E17").Value = '  +2.69%  '
$ws.Range("D18").Value = '3.183.00'
$ws.Range("E18").Value = '  +3.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.68'
$ws.Range("E19").Value = '  -3.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.11'
$ws.Range("E20").Value = '  +9.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.81'
$ws.Range("E21").Value = '  +6.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000203'
$ws.Range("E22").Value = '  -3.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '438.41'
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.04'
$ws.Range("E24").Value = '  +3.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.72'
$ws.Range("E25").Value = '  +2.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.55'
$ws.Range("E26").Value = '  +2.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.82'
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.126'
$ws.Range("E30").Value = '  +45.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.228'
$ws.Range("E31").Value = '  +18.71%  '
$ws.Range("E32").Value = '  +7.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.28'
$ws.Range("E33").Value = '  +1.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("E35").Value = '  +11.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.72'
$ws.Range("E36").Value = '  +9.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.06'
$ws.Range("E37").Value = '  +1.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '502.36'
$ws.Range("E38").Value = '  +0.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.94'
$ws.Range("E39").Value = '  +3.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.34'
$ws.Range("E40").Value = '  +6.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.446'
$ws.Range("E41").Value = '  +11.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.76'
$ws.Range("E42").Value = '  +4.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.42'
$ws.Range("E43").Value = '  -7.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.09'
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.710'
$ws.Range("E46").Value = '  +5.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '158.42'
$ws.Range("E47").Value = '  +4.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.91'
$ws.Range("E48").Value = '  +2.43%  '
$ws.Range("E49").Value = '  +4.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.93'
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.40'
$ws.Range("E51").Value = '  +1.23%  '
